$wb = $excel.ActiveWorkbook

# 1. Text change: every cell showing the status "Ready for handoff" becomes "In Translation"
#    - Overview sheet: columns E (zh-cn) & F (de-de), rows 2-3
#    - zh-cn / de-de sheets: column C (Status), rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# 2. Narrow the "Status" columns (report was regenerated with a tighter auto-fit width)
#    - Overview: columns E & F
#    - zh-cn / de-de: column C
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
